$d = $word.ActiveDocument

# --- A. "Bias" -> "Dark" in the caption run right after the first image ---
# (paragraph: "Bias Frames are an important part of the post-processing ...")
$p = $d.Paragraphs.Item(7)
$p.Range.Find.Execute("Bias", $true, $false, $false, $false, $false, $true, 1, $false, "Dark", 1) | Out-Null

# --- B. Heading: "Precalibrating Input Flats" -> "Precalibrating Input Darks" ---
$p = $d.Paragraphs.Item(13)
$p.Range.Find.Execute("Precalibrating Input Flats", $true, $false, $false, $false, $false, $true, 1, $false, "Precalibrating Input Darks", 1) | Out-Null

# --- E/F. "bias" -> "dark" for the two occurrences inside the worked example paragraph ---
$p = $d.Paragraphs.Item(21)
$p.Range.Find.Execute("5 bias frames for calibration", $true, $false, $false, $false, $false, $true, 1, $false, "5 dark frames for calibration", 1) | Out-Null
$p = $d.Paragraphs.Item(21)
$p.Range.Find.Execute("collection of bias frames as", $true, $false, $false, $false, $false, $true, 1, $false, "collection of dark frames as", 1) | Out-Null

# --- G. " flat frames)." -> " dark frames)." at the end of the same paragraph ---
$p = $d.Paragraphs.Item(21)
$p.Range.Find.Execute("5 flat frames).", $true, $false, $false, $false, $false, $true, 1, $false, "5 dark frames).", 1) | Out-Null

# --- H. "input flat frames" -> "input dark frames" ---
$p = $d.Paragraphs.Item(76)
$p.Range.Find.Execute("input flat frames", $true, $false, $false, $false, $false, $true, 1, $false, "input dark frames", 1) | Out-Null

# --- I. "Flat Frame files" -> "Dark Frame files" ---
$p = $d.Paragraphs.Item(78)
$p.Range.Find.Execute("Flat Frame files", $true, $false, $false, $false, $false, $true, 1, $false, "Dark Frame files", 1) | Out-Null

# --- J. "flat frame files" -> "dark frame files" ---
$p = $d.Paragraphs.Item(79)
$p.Range.Find.Execute("flat frame files", $true, $false, $false, $false, $false, $true, 1, $false, "dark frame files", 1) | Out-Null

# --- K. "not all flat frames" -> "not all dark frames" ---
$p = $d.Paragraphs.Item(85)
$p.Range.Find.Execute("not all flat frames", $true, $false, $false, $false, $false, $true, 1, $false, "not all dark frames", 1) | Out-Null

# --- D. "flat images" -> "dark images" (still in paragraph 18, before paragraph 17 is removed) ---
$p = $d.Paragraphs.Item(18)
$p.Range.Find.Execute("flat images", $true, $false, $false, $false, $false, $true, 1, $false, "dark images", 1) | Out-Null

# --- C. Delete the whole paragraph 17 ("So, ideally you would prepare dark frames for scaling ... combining them.") ---
$d.Paragraphs.Item(17).Range.Delete() | Out-Null
